$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.611.06'
$ws.Range("E2").Value = '  +2.96%  '
$ws.Range("D3").Value = '3.764.05'
$ws.Range("E3").Value = '  +7.29%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '418.81'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.66'
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("D7").Value = '3.747.13'
$ws.Range("E7").Value = '  +7.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.646'
$ws.Range("E8").Value = '  -1.69%  '
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.767'
$ws.Range("E10").Value = '  -2.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.185'
$ws.Range("E11").Value = '  +11.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000403'
$ws.Range("E12").Value = '  +48.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.39'
$ws.Range("E13").Value = '  -2.35%  '
$ws.Range("E14").Value = '  +4.38%  '
$ws.Range("D15").Value = '4.365.17'
$ws.Range("E15").Value = '  +7.62%  '
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").Value = '3.752.27'
$ws.Range("E17").Value = '  +6.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.35'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.15'
$ws.Range("E19").Value = '  +2.47%  '
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("D21").Value = '67.557.54'
$ws.Range("E21").Value = '  +3.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '441.12'
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.52'
$ws.Range("E23").Value = '  +18.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.75'
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.07'
$ws.Range("E25").Value = '  -5.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '38.06'
$ws.Range("E26").Value = '  +11.59%  '
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.00'
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("E29").Value = '  +5.16%  '
$ws.Range("E30").Value = '  +5.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.52'
$ws.Range("E31").Value = '  +0.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.76'
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("E33").Value = '  -2.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.162'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '40.90'
$ws.Range("E35").Value = '  +3.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.02'
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  -3.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.98'
$ws.Range("E39").Value = '  +28.36%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.147'
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '0.0₃0698'
$ws.Range("E41").Value = '  -5.09%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '28.25'
$ws.Range("E42").Value = '  +31.36%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.37'
$ws.Range("E44").Value = '  +3.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.07'
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.20'
$ws.Range("E46").Value = '  +24.94%  '
$ws.Range("E47").Value = '  +4.58%  '
$ws.Range("E48").Value = '  -4.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.27'
$ws.Range("E49").Value = '  -4.86%  '
$ws.Range("E50").Value = '  -7.12%  '
$ws.Range("E51").Value = '  -2.89%  '
